# Insert a new column B ("n" = sample size per guild) ahead of the existing
# percentage columns (which shift from B:D to C:E), then select C10 to match
# the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing B:D columns to C:E by inserting a new column at B.
$ws.Columns.Item(2).Insert()

# New column header + values ("n" = count of species per guild).
$ws.Range("B1").Value = "n"
$ws.Range("B2").Value = 57
$ws.Range("B3").Value = 184
$ws.Range("B4").Value = 23
$ws.Range("B5").Value = 57
$ws.Range("B6").Value = 120
$ws.Range("B7").Value = 38

# Match the saved selection/active cell in the target workbook.
$ws.Range("C10").Select()
